$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one "new report date" worth of price rows at the very top
# of the data block (row 377 onward). A fresh pair of rows (1a/2a nueva(o))
# was published, so the two new rows are inserted at the top and everything
# that was there before shifts down by two rows - growing the used range
# from A1:R466 to A1:R468.
$ws.Rows("377:378").Insert()

# Row 377: "1a nueva(o)"
$ws.Cells.Item(377, 1).Value = 8
$ws.Cells.Item(377, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(377, 3).Value = "Coquimbo"
$ws.Cells.Item(377, 4).Value = 44511
$ws.Cells.Item(377, 5).Value = 4
$ws.Cells.Item(377, 6).Value = 100112004
$ws.Cells.Item(377, 7).Value = "Cebolla"
$ws.Cells.Item(377, 8).Value = "Sin especificar"
$ws.Cells.Item(377, 9).Value = "1a nueva(o)"
$ws.Cells.Item(377, 10).Value = 2600
$ws.Cells.Item(377, 11).Value = 4800
$ws.Cells.Item(377, 12).Value = 5000
$ws.Cells.Item(377, 13).Value = 4900
$ws.Cells.Item(377, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(377, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(377, 16).Value = 272
$ws.Cells.Item(377, 17).Value = 18
$ws.Cells.Item(377, 18).Value = "Hortaliza"

# Row 378: "2a nueva(o)"
$ws.Cells.Item(378, 1).Value = 8
$ws.Cells.Item(378, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(378, 3).Value = "Coquimbo"
$ws.Cells.Item(378, 4).Value = 44511
$ws.Cells.Item(378, 5).Value = 4
$ws.Cells.Item(378, 6).Value = 100112004
$ws.Cells.Item(378, 7).Value = "Cebolla"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "2a nueva(o)"
$ws.Cells.Item(378, 10).Value = 1560
$ws.Cells.Item(378, 11).Value = 4500
$ws.Cells.Item(378, 12).Value = 4600
$ws.Cells.Item(378, 13).Value = 4550
$ws.Cells.Item(378, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(378, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(378, 16).Value = 253
$ws.Cells.Item(378, 17).Value = 18
$ws.Cells.Item(378, 18).Value = "Hortaliza"
